$wb = $excel.ActiveWorkbook

# --- Add the new worksheet after the existing one ---------------------------
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "EquipmentRarityProgression"

# --- Header row ---------------------------------------------------------
# -4108 == xlCenter
$ws.Range("A1:G1").HorizontalAlignment = -4108
$ws.Range("A1:G1").Font.Bold = $true

$ws.Range("A1").Value = "Area"
$ws.Range("B1").Value = "Common"
$ws.Range("C1").Value = "Uncommon"
$ws.Range("D1").Value = "Rare"
$ws.Range("E1").Value = "Epic"
$ws.Range("F1").Value = "Legendary"
$ws.Range("G1").Value = "SUM CHECK"

# --- Data rows -----------------------------------------------------------
$ws.Range("A2:F10").HorizontalAlignment = -4108

$data = @(
  @(1, 75, 20, 5,  0,  0),
  @(2, 50, 40, 10, 0,  0),
  @(3, 20, 45, 30, 5,  0),
  @(4, 20, 45, 30, 5,  0),
  @(5, 20, 30, 40, 10, 0),
  @(6, 20, 30, 40, 10, 0),
  @(7, 15, 20, 50, 10, 5),
  @(8, 10, 30, 35, 10, 15),
  @(9, 5,  10, 35, 20, 30)
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $r = $r + 1
}

# --- "SUM CHECK" column ---------------------------------------------------
$ws.Range("G2").Formula  = "=SUM(B2:F2)"
$ws.Range("G3:G6").Formula = "=SUM(B3:F3)"

# Rows 7-10 were edited by hand afterwards (row 7 still references row 8 -
# an authoring quirk in the source data that we reproduce faithfully).
$ws.Range("G7").Formula  = "=SUM(B8:F8)"
$ws.Range("G8").Formula  = "=SUM(B8:F8)"
$ws.Range("G9").Formula  = "=SUM(B9:F9)"
$ws.Range("G10").Formula = "=SUM(B10:F10)"

# --- Selection / active sheet --------------------------------------------
[void]$ws.Range("C13").Select()

Write-Output "EquipmentRarityProgression sheet added"
